$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 19
$ws1.Range("F3").Value = 175
$ws1.Range("F6").Value = 547
$ws1.Range("F7").Value = 1689
$ws1.Range("F11").Value = 1669
$ws1.Range("F12").Value = 128
$ws1.Range("F13").Value = 82
$ws1.Range("F14").Value = 411
$ws1.Range("F15").Value = 267
$ws1.Range("F16").Value = 195
$ws1.Range("F18").Value = 27
$ws1.Range("F19").Value = 34
$ws1.Range("F21").Value = 270
$ws1.Range("F22").Value = 301
$ws1.Range("F23").Value = 162
$ws1.Range("F24").Value = 228
$ws1.Range("F25").Value = 244

# Sheet "全部类型" (sheet4) - column F ("想去人数") updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 19
$ws4.Range("F3").Value = 175
$ws4.Range("F6").Value = 547
$ws4.Range("F7").Value = 1689
$ws4.Range("F12").Value = 1669
$ws4.Range("F13").Value = 128
$ws4.Range("F14").Value = 82
$ws4.Range("F15").Value = 411
$ws4.Range("F16").Value = 267
$ws4.Range("F17").Value = 195
$ws4.Range("F19").Value = 27
$ws4.Range("F20").Value = 34
$ws4.Range("F22").Value = 270
$ws4.Range("F23").Value = 301
$ws4.Range("F24").Value = 162
$ws4.Range("F25").Value = 228
$ws4.Range("F26").Value = 244
